$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {
    foreach ($hdr in $sec.Headers) {
        if ($hdr.Exists) {
            foreach ($shp in $hdr.Range.InlineShapes) {
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image2.jpg"
                }
            }
        }
    }
    foreach ($ftr in $sec.Footers) {
        if ($ftr.Exists) {
            foreach ($shp in $ftr.Range.InlineShapes) {
                if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Name = "image1.png"
                }
            }
        }
    }
}
